$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Col4a5", "Cd93", "ECs", 1.0, 0.3333333333333333, 0.02193933333333333, 0.065818, 0.00461943947553802, 0.004619439475538019, 3.0, 1.0, 127.984071, 383.952213, 0.4594739502473105, 0.4594739502473105, 2.807885195026, 25.270966755234, 0.002122512103753818, 0.002122512103753818),
    @("FAPs", "Col4a5", "Cd93", "M1", 1.0, 0.3333333333333333, 0.02193933333333333, 0.065818, 0.00461943947553802, 0.004619439475538019, 3.0, 1.0, 60.45343933333334, 181.360318, 0.2170331070069088, 0.2170331070069088, 1.326308156680445, 11.936773410124, 0.001002571302006382, 0.001002571302006381),
    @("FAPs", "Col4a5", "Cd93", "M2", 1.0, 0.3333333333333333, 0.02193933333333333, 0.065818, 0.00461943947553802, 0.004619439475538019, 3.0, 1.0, 59.37981533333334, 178.139446, 0.2131787034353868, 0.2131787034353868, 1.302753561869778, 11.724782056828, 0.0009847661179934384, 0.0009847661179934381),
    @("FAPs", "Col4a5", "Cd93", "Neutro", 1.0, 0.3333333333333333, 0.02193933333333333, 0.065818, 0.00461943947553802, 0.004619439475538019, 3.0, 1.0, 28.25780433333334, 84.773413, 0.1014479761497213, 0.1014479761497213, 0.6199573885371112, 5.579616496834, 0.0004686327857394621, 0.000468632785739462),
    @("FAPs", "Col4a5", "Cd93", "sCs", 1.0, 0.3333333333333333, 0.02193933333333333, 0.065818, 0.00461943947553802, 0.004619439475538019, 3.0, 1.0, 2.469651333333333, 7.408954, 0.008866263160672582, 0.008866263160672582, 0.05418250381911111, 0.487642534372, 0.00004095716604491942, 0.00004095716604491941),
    @("sCs", "Col4a5", "Cd93", "ECs", 3.0, 1.0, 3.059667666666666, 9.179003, 0.6442287642329135, 0.6442287642329135, 3.0, 1.0, 127.984071, 383.952213, 0.4594739502473105, 0.4594739502473105, 391.588723887071, 3524.298514983639, 0.29600633516504, 0.29600633516504),
    @("sCs", "Col4a5", "Cd93", "M1", 3.0, 1.0, 3.059667666666666, 9.179003, 0.6442287642329135, 0.6442287642329135, 3.0, 1.0, 60.45343933333334, 181.360318, 0.2170331070069088, 0.2170331070069088, 184.9674336669949, 1664.706903002954, 0.1398189703246905, 0.1398189703246905),
    @("sCs", "Col4a5", "Cd93", "M2", 3.0, 1.0, 3.059667666666666, 9.179003, 0.6442287642329135, 0.6442287642329135, 3.0, 1.0, 59.37981533333334, 178.139446, 0.2131787034353868, 0.2131787034353868, 181.6825010280376, 1635.142509252338, 0.137335852674954, 0.137335852674954),
    @("sCs", "Col4a5", "Cd93", "Neutro", 3.0, 1.0, 3.059667666666666, 9.179003, 0.6442287642329135, 0.6442287642329135, 3.0, 1.0, 28.25780433333334, 84.773413, 0.1014479761497213, 0.1014479761497213, 86.45949024969322, 778.1354122472391, 0.06535570430886505, 0.06535570430886505),
    @("sCs", "Col4a5", "Cd93", "sCs", 3.0, 1.0, 3.059667666666666, 9.179003, 0.6442287642329135, 0.6442287642329135, 3.0, 1.0, 2.469651333333333, 7.408954, 0.008866263160672582, 0.008866263160672582, 7.556312332540221, 68.006810992862, 0.005711901759363903, 0.005711901759363903),
    @("ECs", "Col4a5", "Cd93", "ECs", 3.0, 1.0, 1.667742666666667, 5.003228, 0.3511517962915484, 0.3511517962915484, 3.0, 1.0, 127.984071, 383.952213, 0.4594739502473105, 0.4594739502473105, 213.444495860396, 1921.000462743564, 0.1613451029785166, 0.1613451029785166),
    @("ECs", "Col4a5", "Cd93", "M1", 3.0, 1.0, 1.667742666666667, 5.003228, 0.3511517962915484, 0.3511517962915484, 3.0, 1.0, 60.45343933333334, 181.360318, 0.2170331070069088, 0.2170331070069088, 100.8207801229449, 907.387021106504, 0.07621156538021186, 0.07621156538021186),
    @("ECs", "Col4a5", "Cd93", "M2", 3.0, 1.0, 1.667742666666667, 5.003228, 0.3511517962915484, 0.3511517962915484, 3.0, 1.0, 59.37981533333334, 178.139446, 0.2131787034353868, 0.2131787034353868, 99.03025157018756, 891.2722641316881, 0.07485808464243937, 0.07485808464243936),
    @("ECs", "Col4a5", "Cd93", "Neutro", 3.0, 1.0, 1.667742666666667, 5.003228, 0.3511517962915484, 0.3511517962915484, 3.0, 1.0, 28.25780433333334, 84.773413, 0.1014479761497213, 0.1014479761497213, 47.12674595301822, 424.140713577164, 0.0356236390551168, 0.0356236390551168),
    @("ECs", "Col4a5", "Cd93", "sCs", 3.0, 1.0, 1.667742666666667, 5.003228, 0.3511517962915484, 0.3511517962915484, 3.0, 1.0, 2.469651333333333, 7.408954, 0.008866263160672582, 0.008866263160672582, 4.118742900390222, 37.068686103512, 0.003113404235263759, 0.003113404235263759),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
